$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 18; $r++) {
    $src = $ws.Range("A" + ($r+1) + ":K" + ($r+1))
    $dst = $ws.Range("A" + $r + ":K" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4104)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
$ws.Range("A19:K19").Clear()
$excel.CutCopyMode = 0

$ws.Range("M4").Select()
